# Updated symbol list (cryptos price/volume snapshot) on Thu Jan 26 13:36:01 UTC 2023
# with GitHub Actions. Updates Price (column D) and Volume(1h) (column E) values
# for various rows. All values are written with a leading apostrophe so Excel
# keeps them as literal text (matching the original inline-string cell contents)
# instead of auto-converting numeric/percent-looking text into real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.29"
$ws.Range("E2").Value = "'1.13%"
$ws.Range("D3").Value = "'35.86"
$ws.Range("E3").Value = "'0.59%"
$ws.Range("D4").Value = "'5.072"
$ws.Range("E4").Value = "'0.39%"
$ws.Range("D5").Value = "'0.08081"
$ws.Range("E5").Value = "'0.49%"
$ws.Range("D6").Value = "'1.956"
$ws.Range("E6").Value = "'1.33%"
$ws.Range("D7").Value = "'4.162"
$ws.Range("E7").Value = "'2.55%"
$ws.Range("D8").Value = "'7.777"
$ws.Range("E8").Value = "'-0.39%"
$ws.Range("D9").Value = "'0.9302"
$ws.Range("E9").Value = "'0.49%"
$ws.Range("D10").Value = "'0.1358"
$ws.Range("E10").Value = "'4.27%"
$ws.Range("D11").Value = "'0.1895"
$ws.Range("E11").Value = "'1.89%"
$ws.Range("D12").Value = "'0.09242"
$ws.Range("E12").Value = "'-1.46%"
$ws.Range("D13").Value = "'0.03535"
$ws.Range("E13").Value = "'3.57%"
$ws.Range("D14").Value = "'0.09870"
$ws.Range("E14").Value = "'0.04%"
$ws.Range("D15").Value = "'0.001409"
$ws.Range("E15").Value = "'0.76%"
$ws.Range("D16").Value = "'0.005811"
$ws.Range("E16").Value = "'0.99%"
$ws.Range("E17").Value = "'1.52%"
$ws.Range("E18").Value = "'-1.11%"
$ws.Range("E19").Value = "'1.20%"
$ws.Range("D20").Value = "'0.1322"
$ws.Range("E20").Value = "'1.60%"
$ws.Range("D21").Value = "'4.887"
$ws.Range("E21").Value = "'-3.06%"
$ws.Range("D22").Value = "'0.2517"
$ws.Range("E22").Value = "'2.02%"
$ws.Range("D23").Value = "'0.04413"
$ws.Range("E23").Value = "'-1.81%"
$ws.Range("D24").Value = "'0.001223"
$ws.Range("E24").Value = "'0.79%"
$ws.Range("D25").Value = "'0.004765"
$ws.Range("E25").Value = "'-1.08%"
$ws.Range("E26").Value = "'32.00%"
$ws.Range("D39").Value = "'0.01951"
$ws.Range("E39").Value = "'1.32%"
$ws.Range("D40").Value = "'0.04946"
$ws.Range("E40").Value = "'3.94%"
$ws.Range("D41").Value = "'0.01084"
$ws.Range("E41").Value = "'12.72%"
$ws.Range("D42").Value = "'0.007650"
$ws.Range("E42").Value = "'4.01%"
$ws.Range("E43").Value = "'3.22%"
$ws.Range("D44").Value = "'0.002101"
$ws.Range("E44").Value = "'-0.46%"
$ws.Range("D45").Value = "'0.01077"
$ws.Range("E45").Value = "'-0.84%"
$ws.Range("D46").Value = "'0.00006384"
$ws.Range("E46").Value = "'1.85%"
$ws.Range("E47").Value = "'0.07%"
$ws.Range("D48").Value = "'63.57"
$ws.Range("E48").Value = "'-1.41%"
$ws.Range("E49").Value = "'-20.03%"
$ws.Range("E50").Value = "'0.07%"
$ws.Range("E51").Value = "'0.07%"
